$d = $word.ActiveDocument
$d.Content.Find.Execute("Python, JavaScript, HTML, CSS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Python, SQL, JavaScript, HTML, CSS", 2)
